$p = $ppt.ActivePresentation

$newStyleId = "{5817E9E1-5101-4347-A747-FCE402E6098B}"

foreach ($slideIndex in 14, 15, 16) {
    $s = $p.Slides.Item($slideIndex)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $sh = $s.Shapes.Item($i)
        if ($sh.HasTable) {
            $sh.Table.ApplyStyle($newStyleId)
        }
    }
}
